$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 108. This pushes the
# existing rows 108..215 down to 110..217, which matches the target layout
# exactly (the two brand-new records end up at rows 108-109).
$ws.Rows.Item(108).Insert()
$ws.Rows.Item(108).Insert()

# Populate the two newly-inserted rows (108 and 109 after the shift) with the
# new data record (both rows carry identical values).
$newDate = 44923
$newVolumen = 40
$newPrecioMin = 28000
$newPrecioMax = 28000
$newPrecioProm = 28000
$newUnidad = "$/malla 25 kilos"
$newOrigen = "Provincia de Limarí"
$newPrecioKg = 1120

foreach ($r in 108,109) {
    $ws.Cells.Item($r, 1).Value = 3
    $ws.Cells.Item($r, 2).Value = "Femacal de La Calera"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = $newDate
    $ws.Cells.Item($r, 5).Value = 5
    $ws.Cells.Item($r, 6).Value = 100112030
    $ws.Cells.Item($r, 7).Value = "Poroto granado"
    $ws.Cells.Item($r, 8).Value = "Sin especificar"
    $ws.Cells.Item($r, 9).Value = "Primera"
    $ws.Cells.Item($r, 10).Value = $newVolumen
    $ws.Cells.Item($r, 11).Value = $newPrecioMin
    $ws.Cells.Item($r, 12).Value = $newPrecioMax
    $ws.Cells.Item($r, 13).Value = $newPrecioProm
    $ws.Cells.Item($r, 14).Value = $newUnidad
    $ws.Cells.Item($r, 15).Value = $newOrigen
    $ws.Cells.Item($r, 16).Value = $newPrecioKg
    $ws.Cells.Item($r, 17).Value = 25
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}
